$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.289.75"
$ws.Range("D3").Value = "3.565.48"
$ws.Range("E3").Value = "  +0.97%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "'606.53"
$ws.Range("E5").Value = "  -0.05%  "
$ws.Range("D6").Value = "'144.70"
$ws.Range("E6").Value = "  +0.17%  "
$ws.Range("D7").Value = "3.564.41"
$ws.Range("E7").Value = "  +1.06%  "
$ws.Range("E8").Value = "  +0.17%  "
$ws.Range("E9").Value = "  +2.02%  "
$ws.Range("E10").Value = "  -0.29%  "
$ws.Range("D11").Value = "'7.80"
$ws.Range("E11").Value = "  -2.90%  "
$ws.Range("E12").Value = "  +0.15%  "
$ws.Range("D13").Value = "4.171.01"
$ws.Range("E13").Value = "  +0.93%  "
$ws.Range("E14").Value = "  -0.28%  "
$ws.Range("D15").Value = "'30.35"
$ws.Range("E15").Value = "  -0.36%  "
$ws.Range("D16").Value = "3.579.39"
$ws.Range("E16").Value = "  +1.25%  "
$ws.Range("D17").Value = "66.334.31"
$ws.Range("E17").Value = "  +0.13%  "
$ws.Range("E18").Value = "  -0.57%  "
$ws.Range("D19").Value = "'11.49"
$ws.Range("E19").Value = "  +4.30%  "
$ws.Range("E20").Value = "  -0.05%  "
$ws.Range("D21").Value = "'14.84"
$ws.Range("E21").Value = "  -1.50%  "
$ws.Range("D22").Value = "'431.51"
$ws.Range("E22").Value = "  +1.35%  "
$ws.Range("E23").Value = "  +1.48%  "
$ws.Range("D24").Value = "'79.59"
$ws.Range("E24").Value = "  +0.86%  "
$ws.Range("D25").Value = "3.707.61"
$ws.Range("E25").Value = "  +0.83%  "
$ws.Range("E26").Value = "  +0.04%  "
$ws.Range("E27").Value = "  -0.20%  "
$ws.Range("D28").Value = "'2.50"
$ws.Range("E28").Value = "  +1.05%  "
$ws.Range("E29").Value = "  -1.54%  "
$ws.Range("D30").Value = "'7.96"
$ws.Range("E30").Value = "  -1.26%  "
$ws.Range("E31").Value = "  -0.06%  "
$ws.Range("E32").Value = "  -2.87%  "
$ws.Range("D33").Value = "3.559.66"
$ws.Range("E33").Value = "  +1.01%  "
$ws.Range("D34").Value = "'25.47"
$ws.Range("E34").Value = "  +0.48%  "
$ws.Range("E35").Value = "  -5.68%  "
$ws.Range("D36").Value = "'7.85"
$ws.Range("E36").Value = "  -0.32%  "
$ws.Range("E38").Value = "  -1.35%  "
$ws.Range("D39").Value = "'5.61"
$ws.Range("E39").Value = "  -0.13%  "
$ws.Range("D40").Value = "'175.76"
$ws.Range("E40").Value = "  +2.63%  "
$ws.Range("D41").Value = "'0.0849"
$ws.Range("E41").Value = "  -1.51%  "
$ws.Range("E42").Value = "  +0.49%  "
$ws.Range("E43").Value = "  -0.17%  "
$ws.Range("D44").Value = "'1.93"
$ws.Range("E44").Value = "  +2.51%  "
$ws.Range("D45").Value = "'46.00"
$ws.Range("E45").Value = "  +1.55%  "
$ws.Range("E46").Value = "  +0.04%  "
$ws.Range("E47").Value = "  +3.98%  "
$ws.Range("E48").Value = "  -1.47%  "
$ws.Range("D49").Value = "'25.12"
$ws.Range("E49").Value = "  -3.33%  "
$ws.Range("E50").Value = "  -0.65%  "
$ws.Range("D51").Value = "'23.39"
$ws.Range("E51").Value = "  +4.13%  "
